$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 151.78572
$ws.Cells.Item(5, 9).Value = 153.125
$ws.Cells.Item(5, 11).Value = 153.125
$ws.Cells.Item(5, 13).Value = -38.125
$ws.Cells.Item(11, 8).Value = 46.117645
$ws.Cells.Item(11, 9).Value = 46.117645
$ws.Cells.Item(11, 11).Value = 46.117645
$ws.Cells.Item(11, 13).Value = 93.88235499999999
$ws.Cells.Item(19, 8).Value = 605.61536
$ws.Cells.Item(19, 9).Value = 495.5
$ws.Cells.Item(19, 10).Value = 700
$ws.Cells.Item(19, 11).Value = 495.5
$ws.Cells.Item(19, 12).Value = 700
$ws.Cells.Item(19, 13).Value = -320.5
$ws.Cells.Item(19, 14).Value = -1050
$ws.Cells.Item(43, 8).Value = 2979
$ws.Cells.Item(43, 9).Value = 2979
$ws.Cells.Item(43, 11).Value = 2979
$ws.Cells.Item(43, 13).Value = -2910
$ws.Cells.Item(105, 8).Value = 29818
$ws.Cells.Item(105, 10).Value = 29818
$ws.Cells.Item(105, 12).Value = 29818
$ws.Cells.Item(105, 14).Value = -36806
$ws.Cells.Item(106, 8).Value = 1444.3334
$ws.Cells.Item(106, 9).Value = 1444.3334
$ws.Cells.Item(106, 11).Value = 1444.3334
$ws.Cells.Item(106, 13).Value = -813.3334
$ws.Cells.Item(132, 8).Value = 4035.9033
$ws.Cells.Item(132, 9).Value = 2107.0833
$ws.Cells.Item(132, 10).Value = 10649
$ws.Cells.Item(132, 11).Value = 6321.249899999999
$ws.Cells.Item(132, 12).Value = 31947
$ws.Cells.Item(132, 13).Value = -3791.249899999999
$ws.Cells.Item(132, 14).Value = -37007
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 817.25
$ws.Cells.Item(2, 9).Value = 817.25
$ws.Cells.Item(2, 11).Value = 817.25
$ws.Cells.Item(2, 13).Value = -704.25
$ws.Cells.Item(26, 8).Value = 1218.1666
$ws.Cells.Item(26, 10).Value = 2989
$ws.Cells.Item(26, 12).Value = 2989
$ws.Cells.Item(26, 14).Value = -3649
$ws.Cells.Item(34, 8).Value = 4750
$ws.Cells.Item(34, 9).Value = 4125
$ws.Cells.Item(34, 10).Value = 6000
$ws.Cells.Item(34, 11).Value = 4125
$ws.Cells.Item(34, 12).Value = 6000
$ws.Cells.Item(34, 13).Value = -3854
$ws.Cells.Item(34, 14).Value = -6542
$ws.Cells.Item(45, 8).Value = 2435.3
$ws.Cells.Item(45, 9).Value = 1640
$ws.Cells.Item(45, 11).Value = 1640
$ws.Cells.Item(45, 13).Value = -1263
$ws.Cells.Item(74, 8).Value = 2575.32
$ws.Cells.Item(74, 9).Value = 2195.1365
$ws.Cells.Item(74, 11).Value = 2195.1365
$ws.Cells.Item(74, 13).Value = -1321.1365
$ws.Cells.Item(77, 8).Value = 2575.32
$ws.Cells.Item(77, 9).Value = 2195.1365
$ws.Cells.Item(77, 11).Value = 10975.6825
$ws.Cells.Item(77, 13).Value = -6607.682500000001
$ws.Cells.Item(96, 8).Value = 2528141
$ws.Cells.Item(96, 10).Value = 2528141
$ws.Cells.Item(96, 12).Value = 2528141
$ws.Cells.Item(96, 14).Value = -2533633
$ws.Cells.Item(106, 8).Value = 39000
$ws.Cells.Item(106, 10).Value = 39000
$ws.Cells.Item(106, 12).Value = 39000
$ws.Cells.Item(106, 14).Value = -41524
$ws.Cells.Item(110, 8).Value = 66668936
$ws.Cells.Item(110, 9).Value = 111112920
$ws.Cells.Item(110, 10).Value = 2962.5
$ws.Cells.Item(110, 11).Value = 111112920
$ws.Cells.Item(110, 12).Value = 2962.5
$ws.Cells.Item(110, 13).Value = -111110875
$ws.Cells.Item(110, 14).Value = -7052.5
$ws.Cells.Item(116, 8).Value = 817.25
$ws.Cells.Item(116, 9).Value = 817.25
$ws.Cells.Item(116, 11).Value = 817.25
$ws.Cells.Item(116, 13).Value = 1476.75
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2778.8572
$ws.Cells.Item(132, 9).Value = 2817.8
$ws.Cells.Item(132, 10).Value = 2000
$ws.Cells.Item(132, 11).Value = 8453.400000000001
$ws.Cells.Item(132, 12).Value = 6000
$ws.Cells.Item(132, 13).Value = -5923.400000000001
$ws.Cells.Item(132, 14).Value = -11060
$ws.Cells.Item(140, 8).Value = 28200
$ws.Cells.Item(140, 9).Value = 28200
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 28200
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = -23020
$ws.Cells.Item(140, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 817.25
$ws.Cells.Item(3, 9).Value = 817.25
$ws.Cells.Item(3, 11).Value = 817.25
$ws.Cells.Item(3, 13).Value = -703.25
$ws.Cells.Item(107, 8).Value = 29417318
$ws.Cells.Item(107, 9).Value = 71430420
$ws.Cells.Item(107, 11).Value = 71430420
$ws.Cells.Item(107, 13).Value = -71428500
$ws.Cells.Item(134, 8).Value = 2498.4
$ws.Cells.Item(134, 9).Value = 764
$ws.Cells.Item(134, 11).Value = 2292
$ws.Cells.Item(134, 13).Value = 243
$ws.Cells.Item(135, 8).Value = 190999.75
$ws.Cells.Item(135, 10).Value = 190999.75
$ws.Cells.Item(135, 12).Value = 190999.75
$ws.Cells.Item(135, 14).Value = -201139.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1971.3
$ws.Cells.Item(22, 9).Value = 1302.1666
$ws.Cells.Item(22, 11).Value = 1302.1666
$ws.Cells.Item(22, 13).Value = -952.1666
$ws.Cells.Item(33, 8).Value = 3919.8
$ws.Cells.Item(33, 9).Value = 1899.7142
$ws.Cells.Item(33, 10).Value = 8633.333000000001
$ws.Cells.Item(33, 11).Value = 1899.7142
$ws.Cells.Item(33, 12).Value = 8633.333000000001
$ws.Cells.Item(33, 13).Value = -1520.7142
$ws.Cells.Item(33, 14).Value = -9391.333000000001
$ws.Cells.Item(58, 8).Value = 2864.2083
$ws.Cells.Item(58, 9).Value = 1672.75
$ws.Cells.Item(58, 10).Value = 5247.125
$ws.Cells.Item(58, 11).Value = 1672.75
$ws.Cells.Item(58, 12).Value = 5247.125
$ws.Cells.Item(58, 13).Value = -1469.75
$ws.Cells.Item(58, 14).Value = -5653.125
$ws.Cells.Item(92, 8).Value = 29999
$ws.Cells.Item(92, 10).Value = 29999
$ws.Cells.Item(92, 12).Value = 29999
$ws.Cells.Item(92, 14).Value = -34991
$ws.Cells.Item(106, 8).Value = 104999.5
$ws.Cells.Item(106, 10).Value = 104999.5
$ws.Cells.Item(106, 12).Value = 104999.5
$ws.Cells.Item(106, 14).Value = -107523.5
$ws.Cells.Item(132, 8).Value = 1850.122
$ws.Cells.Item(132, 9).Value = 1688.7179
$ws.Cells.Item(132, 11).Value = 5066.153700000001
$ws.Cells.Item(132, 13).Value = -2536.153700000001
$ws.Cells.Item(134, 8).Value = 4979.7144
$ws.Cells.Item(134, 9).Value = 2119.3333
$ws.Cells.Item(134, 11).Value = 6357.999899999999
$ws.Cells.Item(134, 13).Value = -3822.999899999999
$ws.Cells.Item(136, 8).Value = 2864.2083
$ws.Cells.Item(136, 9).Value = 1672.75
$ws.Cells.Item(136, 10).Value = 5247.125
$ws.Cells.Item(136, 11).Value = 5018.25
$ws.Cells.Item(136, 12).Value = 15741.375
$ws.Cells.Item(136, 13).Value = -2468.25
$ws.Cells.Item(136, 14).Value = -20841.375
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(115, 8).Value = 12500
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3004.8
$ws.Cells.Item(122, 9).Value = 1500
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(128, 8).Value = 39997.5
$ws.Cells.Item(128, 10).Value = 39997.5
$ws.Cells.Item(128, 12).Value = 39997.5
$ws.Cells.Item(128, 14).Value = -49957.5
$ws.Cells.Item(132, 8).Value = 34174.656
$ws.Cells.Item(132, 9).Value = 41099.96
$ws.Cells.Item(132, 10).Value = 4165
$ws.Cells.Item(132, 11).Value = 123299.88
$ws.Cells.Item(132, 12).Value = 12495
$ws.Cells.Item(132, 13).Value = -120769.88
$ws.Cells.Item(132, 14).Value = -17555
$ws.Cells.Item(135, 8).Value = 222500
$ws.Cells.Item(135, 10).Value = 222500
$ws.Cells.Item(135, 12).Value = 222500
$ws.Cells.Item(135, 14).Value = -232640
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2352
$ws.Cells.Item(55, 9).Value = 2304
$ws.Cells.Item(55, 10).Value = 2400
$ws.Cells.Item(55, 11).Value = 2304
$ws.Cells.Item(55, 12).Value = 2400
$ws.Cells.Item(55, 13).Value = -2131
$ws.Cells.Item(55, 14).Value = -2746
$ws.Cells.Item(93, 8).Value = 1876.9
$ws.Cells.Item(93, 9).Value = 1876.9
$ws.Cells.Item(93, 11).Value = 1876.9
$ws.Cells.Item(93, 13).Value = -628.9000000000001
$ws.Cells.Item(104, 8).Value = 21533.2
$ws.Cells.Item(104, 10).Value = 21533.2
$ws.Cells.Item(104, 12).Value = 21533.2
$ws.Cells.Item(104, 14).Value = -28521.2
$ws.Cells.Item(106, 8).Value = 8758.625
$ws.Cells.Item(106, 10).Value = 8758.625
$ws.Cells.Item(106, 12).Value = 8758.625
$ws.Cells.Item(106, 14).Value = -11282.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 9466.333000000001
$ws.Cells.Item(74, 9).Value = 9199.5
$ws.Cells.Item(74, 10).Value = 10000
$ws.Cells.Item(74, 11).Value = 9199.5
$ws.Cells.Item(74, 12).Value = 10000
$ws.Cells.Item(74, 13).Value = -8263.5
$ws.Cells.Item(74, 14).Value = -11872
$ws.Cells.Item(77, 8).Value = 9466.333000000001
$ws.Cells.Item(77, 9).Value = 9199.5
$ws.Cells.Item(77, 10).Value = 10000
$ws.Cells.Item(77, 11).Value = 27598.5
$ws.Cells.Item(77, 12).Value = 30000
$ws.Cells.Item(77, 13).Value = -22918.5
$ws.Cells.Item(77, 14).Value = -39360
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(104, 8).Value = 34394
$ws.Cells.Item(104, 10).Value = 34394
$ws.Cells.Item(104, 12).Value = 34394
$ws.Cells.Item(104, 14).Value = -41382
$ws.Cells.Item(136, 8).Value = 3348.75
$ws.Cells.Item(136, 9).Value = 2398.3809
$ws.Cells.Item(136, 10).Value = 5163.091
$ws.Cells.Item(136, 11).Value = 2398.3809
$ws.Cells.Item(136, 12).Value = 15489.273
